# feat: add 2022-Q3 data
#
# - Insert a brand-new worksheet named "2022-Q3" between "总计" and the
#   existing "2022-Q2" sheet, populated with the new quarter's fund table
#   (formatted the same way the "总计" sheet's header/index cells are).
# - Update the "总计" (totals) sheet: row 2 becomes the 2022-Q3 summary
#   and a new row 3 is appended holding the 2022-Q2 summary that used to
#   live in row 2 (re-using row 2's styling for the new row's index cell).
# - The pre-existing "2022-Q2" worksheet (with its original fund table)
#   is left untouched.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Give the new row 3 the same index-cell styling as row 2 (bold/bordered)
# before writing into it.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial($xlPasteFormats)

# Row 2 (was 2022-Q2 / 2.15) -> becomes the new 2022-Q3 summary row
$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q3"
$totals.Cells.Item(2, 3).Value = 3
$totals.Cells.Item(2, 4).Value = 1.63

# Row 3 (new) -> carries forward the old 2022-Q2 summary row
$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(3, 2).Value = "2022-Q2"
$totals.Cells.Item(3, 3).Value = 3
$totals.Cells.Item(3, 4).Value = 2.15

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet, positioned right after "总计"
#    (i.e. before the existing "2022-Q2" sheet)
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

# Match formatting with the rest of the workbook: header row + index
# column re-use the same bold/bordered/centered style already used by
# the "总计" sheet's header row and index column.
$totals.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$totals.Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial($xlPasteFormats)

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# The fund-code column ("000988", …) and the D/E/F/G numeric-looking
# text columns must stay TEXT (matching the source data, which stores
# them as inline strings) rather than being auto-coerced into numbers
# (which would also drop the leading zeros on the fund codes). Flip on
# a text number format just long enough to type the value in as text,
# then restore "Normal" cell style so no stray number format lingers on
# the cell (matching the source, which carries no special style there).
$textCells = @()
foreach ($r in 2..4) {
    foreach ($c in 2, 4, 5, 6, 7) {
        $textCells += $q3.Cells.Item($r, $c)
        $q3.Cells.Item($r, $c).NumberFormat = "@"
    }
}

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "000988"
$q3.Cells.Item(2, 3).Value = "嘉实全球互联网股票-人民币（QDII）"
$q3.Cells.Item(2, 4).Value = "10.62"
$q3.Cells.Item(2, 5).Value = "88.35"
$q3.Cells.Item(2, 6).Value = "5.11"
$q3.Cells.Item(2, 7).Value = "0.5427"
$q3.Cells.Item(2, 8).Value = 5

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "000989"
$q3.Cells.Item(3, 3).Value = "嘉实全球互联网股票-美元现汇（QDII）"
$q3.Cells.Item(3, 4).Value = "10.62"
$q3.Cells.Item(3, 5).Value = "88.35"
$q3.Cells.Item(3, 6).Value = "5.11"
$q3.Cells.Item(3, 7).Value = "0.5427"
$q3.Cells.Item(3, 8).Value = 5

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "000990"
$q3.Cells.Item(4, 3).Value = "嘉实全球互联网股票-美元现钞（QDII）"
$q3.Cells.Item(4, 4).Value = "10.62"
$q3.Cells.Item(4, 5).Value = "88.35"
$q3.Cells.Item(4, 6).Value = "5.11"
$q3.Cells.Item(4, 7).Value = "0.5427"
$q3.Cells.Item(4, 8).Value = 5

# Drop the temporary "@" text format again - the source cells carry no
# explicit style/number format at all.
foreach ($cell in $textCells) {
    $cell.Style = "Normal"
}
